# excel data driven test modified
# Rewrites the "Users" sheet from a 6-column (ID/Name/Username/Email/Gender/MobileNo)
# 2-record table into an 11-column (FirstName/LastName/DateOfBirth/Street/PostalCode/
# City/State/Country/Phone/Email/Password) 1-record table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Propagate formatting BEFORE touching values/content, so style slots get
#    reused/deduped against the existing palette instead of minting new ones.
# ---------------------------------------------------------------------------

# Row 1 header style (bordered / centered / wrapped text) -> new columns G1:K1
$ws.Range("F1").Copy()
$ws.Range("G1:K1").PasteSpecial(-4122)

# Row 2 numeric-cell style, grabbed from A2 (ID) BEFORE it gets reformatted to
# text below -> used for the two new numeric columns E2 (PostalCode) and I2
# (Phone), and as the base for the new date column C2 (DateOfBirth).
$ws.Range("A2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C2").NumberFormat = "yyyy-m-d"

# Row 2 text-cell style, grabbed from B2 (Name) -> applied to every row-2 cell
# that should hold plain text in the new layout: A2/F2 (were numeric before:
# ID/MobileNo) plus the brand new G2/H2/J2/K2 columns.
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("J2").PasteSpecial(-4122)
$ws.Range("K2").PasteSpecial(-4122)

# Blank-row style, grabbed from A4 (an already-empty row) -> applied to the
# new columns G:K for every blank row, and to re-blank the old second record
# (row 3) across its original columns A:F.
$ws.Range("A4:F4").Copy()
$ws.Range("G3:K20").PasteSpecial(-4122)
$ws.Range("A3:F3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Header row (row 1): relabel existing columns, add the 5 new ones.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "FirstName"
$ws.Range("B1").Value = "LastName"
$ws.Range("C1").Value = "DateOfBirth"
$ws.Range("D1").Value = "Street"
$ws.Range("E1").Value = "PostalCode"
$ws.Range("F1").Value = "City"
$ws.Range("G1").Value = "State"
$ws.Range("H1").Value = "Country"
$ws.Range("I1").Value = "Phone"
$ws.Range("J1").Value = "Email"
$ws.Range("K1").Value = "Password"

# ---------------------------------------------------------------------------
# 3) Data row (row 2): single Michael Jordan record with the new schema.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Michael"
$ws.Range("B2").Value = "Jordan"
$ws.Range("C2").Value = 30314
$ws.Range("D2").Value = "Bulls Blvd"
$ws.Range("E2").Value = 21998
$ws.Range("F2").Value = "Chicago"
$ws.Range("G2").Value = "IL"
$ws.Range("H2").Value = "United States of America (the)"
$ws.Range("I2").Value = 4079998877
$ws.Range("J2").Value = "mjordan@test.com"
$ws.Range("K2").Value = "Mj1093nciqo!@s9"

$ws.Rows.Item(2).RowHeight = 32.05

# ---------------------------------------------------------------------------
# 4) Drop the now-superfluous last blank row (table shrinks from 21 to 20
#    rows once the second record collapses into a blank row).
# ---------------------------------------------------------------------------
$ws.Rows.Item(21).Delete()

# ---------------------------------------------------------------------------
# 5) Hyperlink: was on D2 (Email, old layout) -> now on J2 (Email, new
#    layout). Re-apply the plain text format afterwards since Hyperlinks.Add
#    auto-switches the cell to the built-in blue/underline Hyperlink style,
#    which the source data does not use (only the string itself is styled).
# ---------------------------------------------------------------------------
$ws.Range("D2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("J2"), "mailto:jordan@test.com", "", "", "mjordan@test.com")
$ws.Range("B2").Copy()
$ws.Range("J2").PasteSpecial(-4122)
